# "break out stock.yaml completed"
# 1) On "10per change": fix up D11:D27 (bsecode) so they are stored as real
#    numbers instead of text, and append 9 new screener rows (28-36).
# 2) On "DND 3 V 0.3": break out the header row (B1:H1) and append the first
#    data row (row 2) that had been missing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "10per change"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("10per change")

$bsecodeFix = @{
    11 = "540699"
    12 = "500510"
    13 = "500300"
    14 = "517354"
    15 = "533274"
    16 = "500493"
    17 = "540180"
    18 = "542066"
    19 = "500112"
    20 = "533096"
    21 = "500440"
    22 = "500093"
    23 = "532843"
    24 = "543396"
    25 = "532898"
    26 = "500477"
    27 = "532155"
}

foreach ($r in $bsecodeFix.Keys) {
    # Assigning a numeric-looking string through .Value lets Excel coerce it
    # to a real number, matching the existing D2:D10 cells.
    $ws1.Cells.Item($r, 4).Value = $bsecodeFix[$r]
}

$newRows = @(
    @{ Row=28; A=1; B="LT";         C="Larsen & Toubro Limited";               D="500510"; E=0.17;  F=3409;     G=10372458; H="05/06/2024 14:17:37" },
    @{ Row=29; A=2; B="LODHA";      C="Macrotech Developers Ltd";              D="543287"; E=-0.53; F=1296.85;  G=1797849;  H="05/06/2024 14:17:37" },
    @{ Row=30; A=3; B="ATGL";       C="Adani Total Gas Ltd";                   D="542066"; E=3.03;  F=936.25;   G=6670432;  H="05/06/2024 14:17:37" },
    @{ Row=31; A=4; B="SBIN";       C="State Bank Of India";                   D="500112"; E=1.88;  F=789.75;   G=74256082; H="05/06/2024 14:17:37" },
    @{ Row=32; A=5; B="ADANIPOWER"; C="Adani Power Limited";                   D="533096"; E=0.51;  F=726.65;   G=34537620; H="05/06/2024 14:17:37" },
    @{ Row=33; A=6; B="CGPOWER";    C="CG Power and Industrial Solutions Ltd"; D="500093"; E=0.1;   F=627.65;   G=10674892; H="05/06/2024 14:17:37" },
    @{ Row=34; A=7; B="PAYTM";      C="One 97 Communications Ltd";             D="543396"; E=-4.91; F=339.85;   G=7404922;  H="05/06/2024 14:17:37" },
    @{ Row=35; A=8; B="POWERGRID";  C="Power Grid Corporation Of India Limited"; D="532898"; E=0.96; F=298.8;   G=45312613; H="05/06/2024 14:17:37" },
    @{ Row=36; A=9; B="GAIL";       C="Gail (india) Limited";                  D="532155"; E=2.55;  F=195.15;   G=49797002; H="05/06/2024 14:17:37" }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws1.Cells.Item($r, 1).Value = $row.A
    $ws1.Cells.Item($r, 2).Value = $row.B
    $ws1.Cells.Item($r, 3).Value = $row.C

    # bsecode must stay text (like the freshly-scraped rows always were)
    $dCell = $ws1.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row.D
    $dCell.ClearFormats()

    $ws1.Cells.Item($r, 5).Value = $row.E
    $ws1.Cells.Item($r, 6).Value = $row.F
    $ws1.Cells.Item($r, 7).Value = $row.G
    $ws1.Cells.Item($r, 8).Value = $row.H
}

# ---------------------------------------------------------------------------
# Sheet 3: "DND 3 V 0.3"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("DND 3 V 0.3")

$ws3.Range("B1").Value = "sr"
$ws3.Range("C1").Value = "nsecode"
$ws3.Range("D1").Value = "name"
$ws3.Range("E1").Value = "bsecode"
$ws3.Range("F1").Value = "per_chg"
$ws3.Range("G1").Value = "close"
$ws3.Range("H1").Value = "volume"

# Match A1's header style (bold + border + centered) on the rest of row 1.
$ws3.Range("A1").Copy()
$ws3.Range("B1:H1").PasteSpecial(-4122)

$ws3.Range("A2").Value = "05/06/2024 14:17:37"
$ws3.Range("B2").Value = 1
$ws3.Range("C2").Value = "PIDILITIND"
$ws3.Range("D2").Value = "Pidilite Industries Limited"

$bseCell = $ws3.Range("E2")
$bseCell.NumberFormat = "@"
$bseCell.Value = "500331"
$bseCell.ClearFormats()

$ws3.Range("F2").Value = 3.4
$ws3.Range("G2").Value = 3166.2
$ws3.Range("H2").Value = 632880

Write-Output "applied chartink_screener break-out edits"
